$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date-serial value that was bumped by one day
# (46081 -> 46082) for every data row (C2:C528).
$ws.Range("C2:C528").Value = 46082
